$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand centers")

# Update the annual demand value for the Luederitz demand center.
$ws.Range("D2").Value = 60638666.7

# Extend the comma-style number formatting one row further down (D7),
# matching the style already used by D2, while leaving it empty.
$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D8").Select()
